# Malta Premier League 2023-2024 sheet update
# - Rows 5 and 6 had their match data (columns F:V) swapped
#   (the "Santa Lucia v Floriana" match moves to row 5, "Birkirkara v Sliema" to row 6;
#   columns A:E are identical between the two rows, so only F:V need to move).
# - Three new match rows are appended at the end of the table (rows 45-47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap rows 5 and 6 (columns F:V) ---
$row5 = @()
for ($c = 6; $c -le 22; $c++) {
  $row5 += ,$ws.Cells.Item(5, $c).Value()
}
$row6 = @()
for ($c = 6; $c -le 22; $c++) {
  $row6 += ,$ws.Cells.Item(6, $c).Value()
}
for ($i = 0; $i -lt 17; $i++) {
  $ws.Cells.Item(5, 6 + $i).Value = $row6[$i]
  $ws.Cells.Item(6, 6 + $i).Value = $row5[$i]
}

# --- 2. Append three new match rows (45, 46, 47) ---
$newRows = @(
  @{ Row=45; A=44; E=45235.45833333334; F='Hamrun';      G=1; H='Mosta';  I=1; J=1.33; K='04/11/2023 02:43'; L=1.4;  M='05/11/2023 09:35'; N=4.6;  O='04/11/2023 02:43'; P=4.66; Q='05/11/2023 09:49'; R=6.84; S='04/11/2023 02:43'; T=7.16; U='05/11/2023 09:49'; V='https://www.betexplorer.com/football/malta/premier-league/hamrun-mosta-fc/v3BHP5e3/' },
  @{ Row=46; A=45; E=45235.58333333334; F='Hibernians';  G=1; H='Gudja';  I=0; J=1.47; K='04/11/2023 02:43'; L=1.61; M='05/11/2023 13:53'; N=3.95; O='04/11/2023 02:43'; P=3.46; Q='05/11/2023 13:53'; R=5.49; S='04/11/2023 02:43'; T=6.15; U='05/11/2023 13:53'; V='https://www.betexplorer.com/football/malta/premier-league/hibernians-gudja/CG9PNqQF/' },
  @{ Row=47; A=46; E=45235.58333333334; F='Santa Lucia'; G=0; H='Sliema'; I=1; J=3.87; K='05/11/2023 12:12'; L=4.4;  M='05/11/2023 13:58'; N=3.34; O='05/11/2023 12:12'; P=3.53; Q='05/11/2023 13:58'; R=1.94; S='05/11/2023 12:12'; T=1.78; U='05/11/2023 13:58'; V='https://www.betexplorer.com/football/malta/premier-league/santa-lucia-sliema/YPALOPA9/' }
)

foreach ($nr in $newRows) {
  $r = $nr.Row

  # Copy the formatting (only) of the last existing data row (44) into the new row,
  # so the index column keeps its bold/border/centered style and the date column
  # keeps its date-time number format - matching the rest of the table exactly.
  $ws.Range("A44:V44").Copy()
  $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)
  $excel.CutCopyMode = $false

  $ws.Cells.Item($r, 1).Value  = $nr.A
  $ws.Cells.Item($r, 2).Value  = "malta"
  $ws.Cells.Item($r, 3).Value  = "premier-league"
  $ws.Cells.Item($r, 4).Value  = "2023-2024"
  $ws.Cells.Item($r, 5).Value  = $nr.E
  $ws.Cells.Item($r, 6).Value  = $nr.F
  $ws.Cells.Item($r, 7).Value  = $nr.G
  $ws.Cells.Item($r, 8).Value  = $nr.H
  $ws.Cells.Item($r, 9).Value  = $nr.I
  $ws.Cells.Item($r, 10).Value = $nr.J
  $ws.Cells.Item($r, 11).Value = $nr.K
  $ws.Cells.Item($r, 12).Value = $nr.L
  $ws.Cells.Item($r, 13).Value = $nr.M
  $ws.Cells.Item($r, 14).Value = $nr.N
  $ws.Cells.Item($r, 15).Value = $nr.O
  $ws.Cells.Item($r, 16).Value = $nr.P
  $ws.Cells.Item($r, 17).Value = $nr.Q
  $ws.Cells.Item($r, 18).Value = $nr.R
  $ws.Cells.Item($r, 19).Value = $nr.S
  $ws.Cells.Item($r, 20).Value = $nr.T
  $ws.Cells.Item($r, 21).Value = $nr.U
  $ws.Cells.Item($r, 22).Value = $nr.V
}
